$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 28.45769926551294
$ws.Range("D2").Value = -0.6923007344870626
$ws.Range("E2").Value = 0.4792803069713264
$ws.Range("C3").Value = 29.07298104513369
$ws.Range("D3").Value = -0.2770189548663069
$ws.Range("E3").Value = 0.07673950135522098
$ws.Range("C4").Value = 28.98648354733675
$ws.Range("D4").Value = -0.3835164526632475
$ws.Range("E4").Value = 0.147084869463401
$ws.Range("C5").Value = 28.90920431921203
$ws.Range("D5").Value = -0.6307956807879727
$ws.Range("E5").Value = 0.3979031909007619
$ws.Range("C6").Value = 30.13803044742587
$ws.Range("D6").Value = 0.5880304474258686
$ws.Range("E6").Value = 0.3457798070998673
$ws.Range("C7").Value = 29.81685577167374
$ws.Range("D7").Value = 0.0668557716737368
$ws.Range("E7").Value = 0.004469694206090827
$ws.Range("C8").Value = 30.75050629839859
$ws.Range("D8").Value = 0.9105062983985874
$ws.Range("E8").Value = 0.8290217194234976
$ws.Range("C9").Value = 30.61713050576026
$ws.Range("D9").Value = 0.8071305057602594
$ws.Range("E9").Value = 0.6514596533288121
$ws.Range("C10").Value = 30.42709721190894
$ws.Range("D10").Value = 0.5070972119089348
$ws.Range("E10").Value = 0.2571475823258151
$ws.Range("C11").Value = 29.69666364499381
$ws.Range("D11").Value = -0.2833363550061918
$ws.Range("E11").Value = 0.08027949006819475
$ws.Range("C12").Value = 29.97555144361638
$ws.Range("D12").Value = -0.06444855638362412
$ws.Range("E12").Value = 0.004153616419933177
$ws.Range("C13").Value = 29.44071962030593
$ws.Range("D13").Value = -0.7692803796940737
$ws.Range("E13").Value = 0.5917923025822582
$ws.Range("C14").Value = 29.86526551103273
$ws.Range("D14").Value = -0.3547344889672672
$ws.Range("E14").Value = 0.1258365576628682
$ws.Range("C15").Value = 29.62894127289632
$ws.Range("D15").Value = -0.7510587271036755
$ws.Range("E15").Value = 0.5640892115585933
$ws.Range("C16").Value = 30.36104736073625
$ws.Range("D16").Value = -0.07895263926375407
$ws.Range("E16").Value = 0.00623351924671248
$ws.Range("C17").Value = 30.39606725906519
$ws.Range("D17").Value = -0.0839327409348094
$ws.Range("E17").Value = 0.00704470500082983
$ws.Range("C18").Value = 30.69627960163829
$ws.Range("D18").Value = 0.006279601638286891
$ws.Range("E18").Value = 0.00003943339673557541
$ws.Range("C19").Value = 30.25991957428342
$ws.Range("D19").Value = -0.4900804257165845
$ws.Range("E19").Value = 0.2401788236705487
$ws.Range("C20").Value = 30.71466799285631
$ws.Range("D20").Value = -0.2253320071436882
$ws.Range("E20").Value = 0.05077451344340314
$ws.Range("C21").Value = 31.03730775975343
$ws.Range("D21").Value = 0.08730775975342908
$ws.Range("E21").Value = 0.007622644913162491
$ws.Range("C22").Value = 31.70770545872092
$ws.Range("D22").Value = 0.6877054587209166
$ws.Range("E22").Value = 0.4729387979545463
$ws.Range("C23").Value = 32.13337620006495
$ws.Range("D23").Value = 1.01337620006495
$ws.Range("E23").Value = 1.026931322858079
$ws.Range("C24").Value = 32.29976812191642
$ws.Range("D24").Value = 1.019768121916421
$ws.Range("E24").Value = 1.039927022476944
$ws.Range("C25").Value = 32.13093560450127
$ws.Range("D25").Value = 0.7509356045012758
$ws.Range("E25").Value = 0.5639042821076966
$ws.Range("C26").Value = 32.54529895307395
$ws.Range("D26").Value = 0.9652989530739546
$ws.Range("E26").Value = 0.9318020688056728
$ws.Range("C27").Value = 32.33996152515173
$ws.Range("D27").Value = 0.6899615251517304
$ws.Range("E27").Value = 0.4760469061897019
$ws.Range("C28").Value = 32.99971959255687
$ws.Range("D28").Value = 1.119719592556873
$ws.Range("E28").Value = 1.253771965955729
$ws.Range("C29").Value = 32.71683441576418
$ws.Range("D29").Value = 0.4368344157641815
$ws.Range("E29").Value = 0.1908243067960338
$ws.Range("C30").Value = 32.81211361359365
$ws.Range("D30").Value = 0.3621136135936425
$ws.Range("E30").Value = 0.1311262691498458
$ws.Range("C31").Value = 33.00473111840309
$ws.Range("D31").Value = 0.1547311184030917
$ws.Range("E31").Value = 0.02394171900227159
$ws.Range("C32").Value = 33.10571218941026
$ws.Range("D32").Value = 0.2057121894102636
$ws.Range("E32").Value = 0.04231750487196417
$ws.Range("C33").Value = 33.0779625320107
$ws.Range("D33").Value = -0.02203746798929984
$ws.Range("E33").Value = 0.0004856499953794152
$ws.Range("C34").Value = 34.06090076721547
$ws.Range("D34").Value = 0.6609007672154732
$ws.Range("E34").Value = 0.4367898241060011
$ws.Range("C35").Value = 33.89679463998552
$ws.Range("D35").Value = 0.1967946399855123
$ws.Range("E35").Value = 0.03872813032702738
$ws.Range("C36").Value = 33.98758366440517
$ws.Range("D36").Value = -0.1124163355948298
$ws.Range("E36").Value = 0.01263743250856939
$ws.Range("C37").Value = 34.27937838195894
$ws.Range("D37").Value = -0.1206216180410564
$ws.Range("E37").Value = 0.01454957473884251
$ws.Range("C38").Value = 34.74813211395009
$ws.Range("D38").Value = -0.1518678860499065
$ws.Range("E38").Value = 0.0230638548132674
$ws.Range("C39").Value = 35.88951697187024
$ws.Range("D39").Value = 0.5895169718702462
$ws.Range("E39").Value = 0.3475302601230646
$ws.Range("C40").Value = 36.26926468172613
$ws.Range("D40").Value = 0.5692646817261249
$ws.Range("E40").Value = 0.3240622778607463
$ws.Range("C41").Value = 36.20843056628006
$ws.Range("D41").Value = -0.09156943371993265
$ws.Range("E41").Value = 0.008384961191789138
$ws.Range("C42").Value = 36.69659092684201
$ws.Range("D42").Value = -0.1034090731579838
$ws.Range("E42").Value = 0.01069343641139324
$ws.Range("C43").Value = 36.68193200156335
$ws.Range("D43").Value = -0.6180679984366506
$ws.Range("E43").Value = 0.3820080506914875
$ws.Range("C44").Value = 37.74507957354775
$ws.Range("D44").Value = -0.1549204264522501
$ws.Range("E44").Value = 0.02400033853214704
$ws.Range("C45").Value = 38.44558104190894
$ws.Range("D45").Value = -0.05441895809105546
$ws.Range("E45").Value = 0.00296142299971605
$ws.Range("C46").Value = 38.8650531928586
$ws.Range("D46").Value = -0.0349468071414023
$ws.Range("E46").Value = 0.001221279329378367
$ws.Range("C47").Value = 39.23220670737122
$ws.Range("D47").Value = -0.1677932926287795
$ws.Range("E47").Value = 0.02815458905120723
$ws.Range("C48").Value = 39.65413175787145
$ws.Range("D48").Value = -0.2458682421285445
$ws.Range("E48").Value = 0.0604511924873806
$ws.Range("C49").Value = 39.70824032249832
$ws.Range("D49").Value = -0.3917596775016818
$ws.Range("E49").Value = 0.1534756449162218
$ws.Range("C50").Value = 39.87678094783252
$ws.Range("D50").Value = -0.7232190521674795
$ws.Range("E50").Value = 0.5230457974180275
$ws.Range("C51").Value = 40.15228882471894
$ws.Range("D51").Value = -0.7477111752810615
$ws.Range("E51").Value = 0.5590720016401862
$ws.Range("C52").Value = 3.570425863113588
$ws.Range("E52").Value = 13.97177902834835
$ws.Range("E53").Value = 0.279435580566967
